$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 709, shifting existing rows 709:805 down to 710:806
$ws.Rows.Item(709).Insert()

# Populate the new row 709 with its data
$ws.Cells.Item(709, 1).Value = 3
$ws.Cells.Item(709, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(709, 3).Value = "Coquimbo"
$ws.Cells.Item(709, 4).Value = 45127
$ws.Cells.Item(709, 5).Value = 5
$ws.Cells.Item(709, 6).Value = 100112032
$ws.Cells.Item(709, 7).Value = "Zapallo italiano"
$ws.Cells.Item(709, 8).Value = "Sin especificar"
$ws.Cells.Item(709, 9).Value = "Primera"
$ws.Cells.Item(709, 10).Value = 105
$ws.Cells.Item(709, 11).Value = 12000
$ws.Cells.Item(709, 12).Value = 13000
$ws.Cells.Item(709, 13).Value = 12524
$ws.Cells.Item(709, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(709, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(709, 16).Value = 209
$ws.Cells.Item(709, 17).Value = 60
$ws.Cells.Item(709, 18).Value = "Hortaliza"
